$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed price (D) / volume change (E) figures for this run.
# Price cells are forced to text number-format before the write so Excel
# does not silently coerce values like "1.00" or "558.30" into numbers
# (which would drop the significant trailing zeros / precision).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.170.76"
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.536.55"
$ws.Range("E3").Value = "  +7.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.45"
$ws.Range("E5").Value = "  +9.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "558.30"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.528.02"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +12.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.98"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +4.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.41"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.093.00"
$ws.Range("E15").Value = "  +7.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.534.00"
$ws.Range("E16").Value = "  +7.69%  "
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.057.54"
$ws.Range("E18").Value = "  +5.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.27"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("E20").Value = "  +8.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.22"
$ws.Range("E22").Value = "  +16.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +8.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.49"
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.15"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.11"
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.35"
$ws.Range("E28").Value = "  +8.82%  "
$ws.Range("E29").Value = "  +10.15%  "
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "640.48"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.58"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.75"
$ws.Range("E33").Value = "  +4.11%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.06"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.48"
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0810"
$ws.Range("E37").Value = "  +9.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.147"
$ws.Range("E38").Value = "  +17.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  +13.44%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.037.45"
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +9.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +10.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0419"
$ws.Range("E47").Value = "  +5.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.60"
$ws.Range("E50").Value = "  +6.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.65"
$ws.Range("E51").Value = "  +9.75%  "
